# Update countries & provincias Spain
# Refresh COVID-19 country statistics and re-sort rows whose totals
# changed enough to swap ranking order, then bump the "last updated" stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos (unchanged rank, just refreshed figures) ---
$ws.Range("B4").Value = 1686442
$ws.Range("C4").Value = 19614
$ws.Range("E4").Value = 1135440
$ws.Range("G4").Value = 617
$ws.Range("H4").Value = 99300

# --- Rows 47-48: Argentina overtakes Dinamarca, rows swap places ---
$ws.Range("A47").Value = "Argentina"
$ws.Range("B47").Value = 12076
$ws.Range("C47").Value = 723
$ws.Range("D47").Value = 3732
$ws.Range("E47").Value = 7892
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 452

$ws.Range("A48").Value = "Dinamarca"
$ws.Range("B48").Value = 11360
$ws.Range("C48").Value = 71
$ws.Range("D48").Value = 9900
$ws.Range("E48").Value = 898
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 562

# --- Rows 212-214: Sahara Occidental overtakes the other two, rows swap ---
$ws.Range("A212").Value = "Sahara Occidental"
$ws.Range("B212").Value = 9
$ws.Range("C212").Value = 3
$ws.Range("D212").Value = 6
$ws.Range("E212").Value = 3
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 6
$ws.Range("E213").Value = 1
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# --- Timestamp footer (title cell, row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 02:35"
